$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.426.46'
$ws.Range("E2").Value = '  -2.89%  '
$ws.Range("D3").Value = '2.486.10'
$ws.Range("E3").Value = '  -1.85%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.92%  '
$ws.Range("E7").Value = '  -2.67%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -3.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.62'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.16%  '
$ws.Range("E11").Value = '  -2.83%  '
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("E13").Value = '  -3.64%  '
$ws.Range("D14").Value = '2.871.67'
$ws.Range("E14").Value = '  -1.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.53'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("D16").Value = '2.484.36'
$ws.Range("E16").Value = '  -2.02%  '
$ws.Range("E17").Value = '  -2.19%  '
$ws.Range("D18").Value = '41.415.86'
$ws.Range("E18").Value = '  -2.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.62%  '
$ws.Range("D20").Value = '0.0₃0925'
$ws.Range("E20").Value = '  -2.43%  '
$ws.Range("E21").Value = '  -8.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.74'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.84%  '
$ws.Range("E24").Value = '  -3.34%  '
$ws.Range("E25").Value = '  -4.77%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.25'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.52%  '
$ws.Range("E29").Value = '  -3.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '152.29'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.86%  '
$ws.Range("E32").Value = '  -6.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.56'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.88%  '
$ws.Range("E34").Value = '  -2.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0751'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.84'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.41%  '
$ws.Range("E37").Value = '  -2.22%  '
$ws.Range("E38").Value = '  -3.86%  '
$ws.Range("E39").Value = '  -2.16%  '
$ws.Range("E40").Value = '  -7.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.21'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.11%  '
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("D43").Value = '2.009.37'
$ws.Range("E43").Value = '  +0.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.52'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -10.65%  '
$ws.Range("E45").Value = '  -3.75%  '
$ws.Range("E46").Value = '  -8.55%  '
$ws.Range("E47").Value = '  -2.94%  '
$ws.Range("D48").Value = '2.734.57'
$ws.Range("E48").Value = '  -1.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '70.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '97.51'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.34%  '
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.48%  '
